$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns A (Company Name), B (Company Number), H (Category)
# for rows 2-8, reflecting the reshuffled company data.
$data = @(
    @{ Row = 2; A = "DGPI LTD"; B = "SC849118"; H = "GP" },
    @{ Row = 3; A = "ST GEORGE CAPITAL (LAND) LIMITED"; B = "16462880"; H = "Capital" },
    @{ Row = 4; A = "AFROSCOT VENTURES LTD"; B = "16462878"; H = "Ventures" },
    @{ Row = 5; A = "DAVIDSON CAPITAL HOLDINGS LTD"; B = "SC849117"; H = "Capital" },
    @{ Row = 6; A = "SAMVIV PARTNERS LTD"; B = "16460672"; H = "Partners" },
    @{ Row = 7; A = "T GILPIN PHYSIO CONSULTANCY LTD"; B = "16460503"; H = "LP" },
    @{ Row = 8; A = "4D CAPITAL PROPCO (44) LIMITED"; B = "16461269"; H = "Capital" }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Range("A$r").Value = $entry.A
    $ws.Range("B$r").Value = $entry.B
    $ws.Range("H$r").Value = $entry.H
}
